$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text while we write the new values so that
# numeric-looking strings (e.g. "0.970") keep their exact original text
# representation instead of being auto-converted to a Double by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.599.51"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "1.581.19"
$ws.Range("E3").Value = "  -3.19%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "206.33"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "22.17"
$ws.Range("E8").Value = "  -5.34%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "1.803.34"
$ws.Range("D13").Value = "1.557.98"
$ws.Range("E13").Value = "  -5.66%  "
$ws.Range("E14").Value = "  -4.51%  "
$ws.Range("D15").Value = "0.529"
$ws.Range("E15").Value = "  -6.48%  "
$ws.Range("D16").Value = "27.547.33"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "62.75"
$ws.Range("D18").Value = "217.54"
$ws.Range("E18").Value = "  -5.03%  "
$ws.Range("D19").Value = "7.35"
$ws.Range("E19").Value = "  -4.47%  "
$ws.Range("D20").Value = "0.0₃0693"
$ws.Range("E20").Value = "  -3.68%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  -4.65%  "
$ws.Range("D23").Value = "9.47"
$ws.Range("E23").Value = "  -5.93%  "
$ws.Range("E24").Value = "  -4.36%  "
$ws.Range("D25").Value = "153.73"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "6.69"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("E28").Value = "  -3.27%  "
$ws.Range("E29").Value = "  -4.68%  "
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("D31").Value = "0.0463"
$ws.Range("E31").Value = "  -3.92%  "
$ws.Range("E32").Value = "  -5.18%  "
$ws.Range("D33").Value = "1.367.92"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  -5.42%  "
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("D36").Value = "0.970"
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.77"
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.17"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "63.52"
$ws.Range("E45").Value = "  -3.62%  "
$ws.Range("D46").Value = "5.25"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("D47").Value = "1.714.07"
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").Value = "87.37"
$ws.Range("D49").Value = "0.0₇0996"
$ws.Range("E49").Value = "  -3.53%  "
$ws.Range("D50").Value = "0.0971"
$ws.Range("E50").Value = "  -4.71%  "
$ws.Range("E51").Value = "  -1.72%  "

# Restore default (General) formatting on column D so the saved file has no
# leftover per-cell number-format/style differences versus the original.
$ws.Range("D2:D51").ClearFormats()
